# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E23) on the account-statement sheet held
# eight consecutive billing periods (1607..1702) in ascending order. This
# edit removes the old periods and enters the new ones (same eight values)
# in descending / most-recent-first order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1702"
$ws.Range("E17").Value = "1701"
$ws.Range("E18").Value = "1612"
$ws.Range("E19").Value = "1611"
$ws.Range("E20").Value = "1610"
$ws.Range("E21").Value = "1609"
$ws.Range("E22").Value = "1608"
$ws.Range("E23").Value = "1607"
